# KP-11725 D: Extension of questionnaire's translation files
#
# Inserts a new "Variable" column (B) into the Translations sheet, right
# after the "Entity Id" column, shifting the existing Type/Index/Original/
# Translation columns one position to the right (B:E -> C:F). The header
# cell gets "Variable" and every data row (2-13) is stamped with the
# variable name "c1".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before the current column B (Type); this shifts the
# old B:E columns (Type, Index, Original, Translation) to C:F and extends
# the used range from A1:E13 to A1:F13.
$ws.Columns("B:B").Insert()

# Header for the newly inserted column.
$ws.Range("B1").Value = "Variable"

# Every translated entity row references questionnaire variable "c1".
$ws.Range("B2:B13").Value = "c1"

# Match the author's final selection in the sheet.
$ws.Range("E6").Select()
